$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("optimization_parameters")
$ws.Range("A1").Value = "test"
